# --- Row 2 (OPEX) : add E2/F2 validity-window strings, move the PERFORM
#     formula from F2 to H2 (now referencing the new B2..F2 layout), and
#     add J2 (new budget-group id) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "9999-12-31 23:59:59+07"
$ws.Range("E2").Value = "2020-01-01 00:00:00+07"
$ws.Range("E2:F2").NumberFormat = "@"
$ws.Range("E2:F2").Font.Name = "Arial Narrow"
$ws.Range("E2:F2").Font.Size = 10

$ws.Range("H2").Formula = '=CONCATENATE("PERFORM ""SchData-OLTP-Budgeting"".""Func_TblBudgetGroup_SET""(varSystemLoginSession, null, null, null, varInstitutionBranchID, ", B2, "::bigint, ''", C2, "''::varchar, ''", D2, "''::varchar, ''", E2, "''::timestamptz, ''", F2, "''::timestamptz);")'

$j2 = $ws.Range("J2")
$j2.Value = 109000000000001
$j2.NumberFormat = "0"
$j2.Font.Name = "Arial Narrow"
$j2.Font.Size = 10
$j2.Interior.Color = 5296274

# --- Row 3 (CAPEX): same treatment; F3 previously held the old literal
#     PERFORM string, it gets replaced by the validity end-date ---
$ws.Range("E3").Value = "2020-01-01 00:00:00+07"
$ws.Range("F3").Value = "9999-12-31 23:59:59+07"
$ws.Range("E3:F3").NumberFormat = "@"
$ws.Range("E3:F3").Font.Name = "Arial Narrow"
$ws.Range("E3:F3").Font.Size = 10

$ws.Range("H3").Formula = '=CONCATENATE("PERFORM ""SchData-OLTP-Budgeting"".""Func_TblBudgetGroup_SET""(varSystemLoginSession, null, null, null, varInstitutionBranchID, ", B3, "::bigint, ''", C3, "''::varchar, ''", D3, "''::varchar, ''", E3, "''::timestamptz, ''", F3, "''::timestamptz);")'

$j3 = $ws.Range("J3")
$j3.Value = 109000000000002
$j3.NumberFormat = "0"
$j3.Font.Name = "Arial Narrow"
$j3.Font.Size = 10
$j3.Interior.Color = 5296274

# --- column widths ---
$ws.Range("E1:F1").ColumnWidth = 18.28515625
$ws.Range("J1").EntireColumn.ColumnWidth = 14

# --- selection / view ---
$ws.Range("G14").Select()

# --- page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "done"
